# Applies the edits described by the commit diff:
#  - Row 17 (Decide on content and website placement):
#      Time Cost (C17): 3 -> 4
#      Time spent (D17): "40mins" -> "1hr20mins"
#  - Row 21 (Design timeline):
#      Time spent (D21): (blank) -> "30mins"
#  - Selection/view moved to reflect where the user was working (C21)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update time cost for "Decide on content and website placement"
$ws.Range("C17").Value = 4

# Update the time-spent note for that same row
$ws.Range("D17").Value = "1hr20mins"

# Add the time-spent note for "Design timeline" row
$ws.Range("D21").Value = "30mins"

# Update view state to match where editing occurred
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("C21").Select()
